$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Random")

# --- Row 1: new block headers (Digits target counts) ---
$ws.Range("AC1").Value = 3500
$ws.Range("AG1").Value = 4000
$ws.Range("AK1").Value = 4500
$ws.Range("AO1").Value = 5000

# --- Row 2: new/backfilled labels ---
$ws.Range("W2").Value = "Num/sec"
$ws.Range("AA2").Value = "Num/sec"
$ws.Range("AE2").Value = "Num/sec"
$ws.Range("AI2").Value = "Num/Sec"
$ws.Range("AM2").Value = "Num/Sec"

# --- Rows 4-33: four new data blocks (AC:AE, AG:AI, AK:AM, AO:AQ) ---
$ws.Range("AC4").Value = 748
$ws.Range("AD4").Formula = "=AC4/1000"
$ws.Range("AE4").Formula = '=$AC$1/AD4'
$ws.Range("AG4").Value = 867
$ws.Range("AH4").Formula = "=AG4/1000"
$ws.Range("AI4").Formula = '=$AG$1/AH4'
$ws.Range("AK4").Value = 976
$ws.Range("AL4").Formula = "=AK4/1000"
$ws.Range("AM4").Formula = '=$AK$1/AL4'
$ws.Range("AO4").Value = 1098
$ws.Range("AP4").Formula = "=AO4/1000"
$ws.Range("AQ4").Formula = '=$AO$1/AP4'

$ws.Range("AC5").Value = 749
$ws.Range("AD5").Formula = "=AC5/1000"
$ws.Range("AE5").Formula = '=$AC$1/AD5'
$ws.Range("AG5").Value = 866
$ws.Range("AH5").Formula = "=AG5/1000"
$ws.Range("AI5").Formula = '=$AG$1/AH5'
$ws.Range("AK5").Value = 1060
$ws.Range("AL5").Formula = "=AK5/1000"
$ws.Range("AM5").Formula = '=$AK$1/AL5'
$ws.Range("AO5").Value = 1180
$ws.Range("AP5").Formula = "=AO5/1000"
$ws.Range("AQ5").Formula = '=$AO$1/AP5'

$ws.Range("AC6").Value = 743
$ws.Range("AD6").Formula = "=AC6/1000"
$ws.Range("AE6").Formula = '=$AC$1/AD6'
$ws.Range("AG6").Value = 862
$ws.Range("AH6").Formula = "=AG6/1000"
$ws.Range("AI6").Formula = '=$AG$1/AH6'
$ws.Range("AK6").Value = 991
$ws.Range("AL6").Formula = "=AK6/1000"
$ws.Range("AM6").Formula = '=$AK$1/AL6'
$ws.Range("AO6").Value = 1141
$ws.Range("AP6").Formula = "=AO6/1000"
$ws.Range("AQ6").Formula = '=$AO$1/AP6'

$ws.Range("AC7").Value = 745
$ws.Range("AD7").Formula = "=AC7/1000"
$ws.Range("AE7").Formula = '=$AC$1/AD7'
$ws.Range("AG7").Value = 864
$ws.Range("AH7").Formula = "=AG7/1000"
$ws.Range("AI7").Formula = '=$AG$1/AH7'
$ws.Range("AK7").Value = 979
$ws.Range("AL7").Formula = "=AK7/1000"
$ws.Range("AM7").Formula = '=$AK$1/AL7'
$ws.Range("AO7").Value = 1091
$ws.Range("AP7").Formula = "=AO7/1000"
$ws.Range("AQ7").Formula = '=$AO$1/AP7'

$ws.Range("AC8").Value = 747
$ws.Range("AD8").Formula = "=AC8/1000"
$ws.Range("AE8").Formula = '=$AC$1/AD8'
$ws.Range("AG8").Value = 863
$ws.Range("AH8").Formula = "=AG8/1000"
$ws.Range("AI8").Formula = '=$AG$1/AH8'
$ws.Range("AK8").Value = 985
$ws.Range("AL8").Formula = "=AK8/1000"
$ws.Range("AM8").Formula = '=$AK$1/AL8'
$ws.Range("AO8").Value = 1081
$ws.Range("AP8").Formula = "=AO8/1000"
$ws.Range("AQ8").Formula = '=$AO$1/AP8'

$ws.Range("AC9").Value = 752
$ws.Range("AD9").Formula = "=AC9/1000"
$ws.Range("AE9").Formula = '=$AC$1/AD9'
$ws.Range("AG9").Value = 860
$ws.Range("AH9").Formula = "=AG9/1000"
$ws.Range("AI9").Formula = '=$AG$1/AH9'
$ws.Range("AK9").Value = 972
$ws.Range("AL9").Formula = "=AK9/1000"
$ws.Range("AM9").Formula = '=$AK$1/AL9'
$ws.Range("AO9").Value = 1098
$ws.Range("AP9").Formula = "=AO9/1000"
$ws.Range("AQ9").Formula = '=$AO$1/AP9'

$ws.Range("AC10").Value = 747
$ws.Range("AD10").Formula = "=AC10/1000"
$ws.Range("AE10").Formula = '=$AC$1/AD10'
$ws.Range("AG10").Value = 855
$ws.Range("AH10").Formula = "=AG10/1000"
$ws.Range("AI10").Formula = '=$AG$1/AH10'
$ws.Range("AK10").Value = 967
$ws.Range("AL10").Formula = "=AK10/1000"
$ws.Range("AM10").Formula = '=$AK$1/AL10'
$ws.Range("AO10").Value = 1075
$ws.Range("AP10").Formula = "=AO10/1000"
$ws.Range("AQ10").Formula = '=$AO$1/AP10'

$ws.Range("AC11").Value = 744
$ws.Range("AD11").Formula = "=AC11/1000"
$ws.Range("AE11").Formula = '=$AC$1/AD11'
$ws.Range("AG11").Value = 862
$ws.Range("AH11").Formula = "=AG11/1000"
$ws.Range("AI11").Formula = '=$AG$1/AH11'
$ws.Range("AK11").Value = 975
$ws.Range("AL11").Formula = "=AK11/1000"
$ws.Range("AM11").Formula = '=$AK$1/AL11'
$ws.Range("AO11").Value = 1075
$ws.Range("AP11").Formula = "=AO11/1000"
$ws.Range("AQ11").Formula = '=$AO$1/AP11'

$ws.Range("AC12").Value = 747
$ws.Range("AD12").Formula = "=AC12/1000"
$ws.Range("AE12").Formula = '=$AC$1/AD12'
$ws.Range("AG12").Value = 866
$ws.Range("AH12").Formula = "=AG12/1000"
$ws.Range("AI12").Formula = '=$AG$1/AH12'
$ws.Range("AK12").Value = 980
$ws.Range("AL12").Formula = "=AK12/1000"
$ws.Range("AM12").Formula = '=$AK$1/AL12'
$ws.Range("AO12").Value = 1085
$ws.Range("AP12").Formula = "=AO12/1000"
$ws.Range("AQ12").Formula = '=$AO$1/AP12'

$ws.Range("AC13").Value = 740
$ws.Range("AD13").Formula = "=AC13/1000"
$ws.Range("AE13").Formula = '=$AC$1/AD13'
$ws.Range("AG13").Value = 866
$ws.Range("AH13").Formula = "=AG13/1000"
$ws.Range("AI13").Formula = '=$AG$1/AH13'
$ws.Range("AK13").Value = 1007
$ws.Range("AL13").Formula = "=AK13/1000"
$ws.Range("AM13").Formula = '=$AK$1/AL13'
$ws.Range("AO13").Value = 1090
$ws.Range("AP13").Formula = "=AO13/1000"
$ws.Range("AQ13").Formula = '=$AO$1/AP13'

$ws.Range("AC14").Value = 745
$ws.Range("AD14").Formula = "=AC14/1000"
$ws.Range("AE14").Formula = '=$AC$1/AD14'
$ws.Range("AG14").Value = 867
$ws.Range("AH14").Formula = "=AG14/1000"
$ws.Range("AI14").Formula = '=$AG$1/AH14'
$ws.Range("AK14").Value = 1002
$ws.Range("AL14").Formula = "=AK14/1000"
$ws.Range("AM14").Formula = '=$AK$1/AL14'
$ws.Range("AO14").Value = 1088
$ws.Range("AP14").Formula = "=AO14/1000"
$ws.Range("AQ14").Formula = '=$AO$1/AP14'

$ws.Range("AC15").Value = 749
$ws.Range("AD15").Formula = "=AC15/1000"
$ws.Range("AE15").Formula = '=$AC$1/AD15'
$ws.Range("AG15").Value = 876
$ws.Range("AH15").Formula = "=AG15/1000"
$ws.Range("AI15").Formula = '=$AG$1/AH15'
$ws.Range("AK15").Value = 977
$ws.Range("AL15").Formula = "=AK15/1000"
$ws.Range("AM15").Formula = '=$AK$1/AL15'
$ws.Range("AO15").Value = 1100
$ws.Range("AP15").Formula = "=AO15/1000"
$ws.Range("AQ15").Formula = '=$AO$1/AP15'

$ws.Range("AC16").Value = 742
$ws.Range("AD16").Formula = "=AC16/1000"
$ws.Range("AE16").Formula = '=$AC$1/AD16'
$ws.Range("AG16").Value = 868
$ws.Range("AH16").Formula = "=AG16/1000"
$ws.Range("AI16").Formula = '=$AG$1/AH16'
$ws.Range("AK16").Value = 990
$ws.Range("AL16").Formula = "=AK16/1000"
$ws.Range("AM16").Formula = '=$AK$1/AL16'
$ws.Range("AO16").Value = 1121
$ws.Range("AP16").Formula = "=AO16/1000"
$ws.Range("AQ16").Formula = '=$AO$1/AP16'

$ws.Range("AC17").Value = 765
$ws.Range("AD17").Formula = "=AC17/1000"
$ws.Range("AE17").Formula = '=$AC$1/AD17'
$ws.Range("AG17").Value = 859
$ws.Range("AH17").Formula = "=AG17/1000"
$ws.Range("AI17").Formula = '=$AG$1/AH17'
$ws.Range("AK17").Value = 975
$ws.Range("AL17").Formula = "=AK17/1000"
$ws.Range("AM17").Formula = '=$AK$1/AL17'
$ws.Range("AO17").Value = 1086
$ws.Range("AP17").Formula = "=AO17/1000"
$ws.Range("AQ17").Formula = '=$AO$1/AP17'

$ws.Range("AC18").Value = 749
$ws.Range("AD18").Formula = "=AC18/1000"
$ws.Range("AE18").Formula = '=$AC$1/AD18'
$ws.Range("AG18").Value = 855
$ws.Range("AH18").Formula = "=AG18/1000"
$ws.Range("AI18").Formula = '=$AG$1/AH18'
$ws.Range("AK18").Value = 967
$ws.Range("AL18").Formula = "=AK18/1000"
$ws.Range("AM18").Formula = '=$AK$1/AL18'
$ws.Range("AO18").Value = 1091
$ws.Range("AP18").Formula = "=AO18/1000"
$ws.Range("AQ18").Formula = '=$AO$1/AP18'

$ws.Range("AC19").Value = 783
$ws.Range("AD19").Formula = "=AC19/1000"
$ws.Range("AE19").Formula = '=$AC$1/AD19'
$ws.Range("AG19").Value = 865
$ws.Range("AH19").Formula = "=AG19/1000"
$ws.Range("AI19").Formula = '=$AG$1/AH19'
$ws.Range("AK19").Value = 986
$ws.Range("AL19").Formula = "=AK19/1000"
$ws.Range("AM19").Formula = '=$AK$1/AL19'
$ws.Range("AO19").Value = 1085
$ws.Range("AP19").Formula = "=AO19/1000"
$ws.Range("AQ19").Formula = '=$AO$1/AP19'

$ws.Range("AC20").Value = 772
$ws.Range("AD20").Formula = "=AC20/1000"
$ws.Range("AE20").Formula = '=$AC$1/AD20'
$ws.Range("AG20").Value = 859
$ws.Range("AH20").Formula = "=AG20/1000"
$ws.Range("AI20").Formula = '=$AG$1/AH20'
$ws.Range("AK20").Value = 987
$ws.Range("AL20").Formula = "=AK20/1000"
$ws.Range("AM20").Formula = '=$AK$1/AL20'
$ws.Range("AO20").Value = 1077
$ws.Range("AP20").Formula = "=AO20/1000"
$ws.Range("AQ20").Formula = '=$AO$1/AP20'

$ws.Range("AC21").Value = 768
$ws.Range("AD21").Formula = "=AC21/1000"
$ws.Range("AE21").Formula = '=$AC$1/AD21'
$ws.Range("AG21").Value = 874
$ws.Range("AH21").Formula = "=AG21/1000"
$ws.Range("AI21").Formula = '=$AG$1/AH21'
$ws.Range("AK21").Value = 986
$ws.Range("AL21").Formula = "=AK21/1000"
$ws.Range("AM21").Formula = '=$AK$1/AL21'
$ws.Range("AO21").Value = 1082
$ws.Range("AP21").Formula = "=AO21/1000"
$ws.Range("AQ21").Formula = '=$AO$1/AP21'

$ws.Range("AC22").Value = 759
$ws.Range("AD22").Formula = "=AC22/1000"
$ws.Range("AE22").Formula = '=$AC$1/AD22'
$ws.Range("AG22").Value = 875
$ws.Range("AH22").Formula = "=AG22/1000"
$ws.Range("AI22").Formula = '=$AG$1/AH22'
$ws.Range("AK22").Value = 971
$ws.Range("AL22").Formula = "=AK22/1000"
$ws.Range("AM22").Formula = '=$AK$1/AL22'
$ws.Range("AO22").Value = 1088
$ws.Range("AP22").Formula = "=AO22/1000"
$ws.Range("AQ22").Formula = '=$AO$1/AP22'

$ws.Range("AC23").Value = 775
$ws.Range("AD23").Formula = "=AC23/1000"
$ws.Range("AE23").Formula = '=$AC$1/AD23'
$ws.Range("AG23").Value = 874
$ws.Range("AH23").Formula = "=AG23/1000"
$ws.Range("AI23").Formula = '=$AG$1/AH23'
$ws.Range("AK23").Value = 978
$ws.Range("AL23").Formula = "=AK23/1000"
$ws.Range("AM23").Formula = '=$AK$1/AL23'
$ws.Range("AO23").Value = 1090
$ws.Range("AP23").Formula = "=AO23/1000"
$ws.Range("AQ23").Formula = '=$AO$1/AP23'

$ws.Range("AC24").Value = 781
$ws.Range("AD24").Formula = "=AC24/1000"
$ws.Range("AE24").Formula = '=$AC$1/AD24'
$ws.Range("AG24").Value = 855
$ws.Range("AH24").Formula = "=AG24/1000"
$ws.Range("AI24").Formula = '=$AG$1/AH24'
$ws.Range("AK24").Value = 974
$ws.Range("AL24").Formula = "=AK24/1000"
$ws.Range("AM24").Formula = '=$AK$1/AL24'
$ws.Range("AO24").Value = 1085
$ws.Range("AP24").Formula = "=AO24/1000"
$ws.Range("AQ24").Formula = '=$AO$1/AP24'

$ws.Range("AC25").Value = 818
$ws.Range("AD25").Formula = "=AC25/1000"
$ws.Range("AE25").Formula = '=$AC$1/AD25'
$ws.Range("AG25").Value = 854
$ws.Range("AH25").Formula = "=AG25/1000"
$ws.Range("AI25").Formula = '=$AG$1/AH25'
$ws.Range("AK25").Value = 963
$ws.Range("AL25").Formula = "=AK25/1000"
$ws.Range("AM25").Formula = '=$AK$1/AL25'
$ws.Range("AO25").Value = 1090
$ws.Range("AP25").Formula = "=AO25/1000"
$ws.Range("AQ25").Formula = '=$AO$1/AP25'

$ws.Range("AC26").Value = 772
$ws.Range("AD26").Formula = "=AC26/1000"
$ws.Range("AE26").Formula = '=$AC$1/AD26'
$ws.Range("AG26").Value = 870
$ws.Range("AH26").Formula = "=AG26/1000"
$ws.Range("AI26").Formula = '=$AG$1/AH26'
$ws.Range("AK26").Value = 985
$ws.Range("AL26").Formula = "=AK26/1000"
$ws.Range("AM26").Formula = '=$AK$1/AL26'
$ws.Range("AO26").Value = 1077
$ws.Range("AP26").Formula = "=AO26/1000"
$ws.Range("AQ26").Formula = '=$AO$1/AP26'

$ws.Range("AC27").Value = 764
$ws.Range("AD27").Formula = "=AC27/1000"
$ws.Range("AE27").Formula = '=$AC$1/AD27'
$ws.Range("AG27").Value = 863
$ws.Range("AH27").Formula = "=AG27/1000"
$ws.Range("AI27").Formula = '=$AG$1/AH27'
$ws.Range("AK27").Value = 979
$ws.Range("AL27").Formula = "=AK27/1000"
$ws.Range("AM27").Formula = '=$AK$1/AL27'
$ws.Range("AO27").Value = 1092
$ws.Range("AP27").Formula = "=AO27/1000"
$ws.Range("AQ27").Formula = '=$AO$1/AP27'

$ws.Range("AC28").Value = 762
$ws.Range("AD28").Formula = "=AC28/1000"
$ws.Range("AE28").Formula = '=$AC$1/AD28'
$ws.Range("AG28").Value = 857
$ws.Range("AH28").Formula = "=AG28/1000"
$ws.Range("AI28").Formula = '=$AG$1/AH28'
$ws.Range("AK28").Value = 969
$ws.Range("AL28").Formula = "=AK28/1000"
$ws.Range("AM28").Formula = '=$AK$1/AL28'
$ws.Range("AO28").Value = 1091
$ws.Range("AP28").Formula = "=AO28/1000"
$ws.Range("AQ28").Formula = '=$AO$1/AP28'

$ws.Range("AC29").Value = 799
$ws.Range("AD29").Formula = "=AC29/1000"
$ws.Range("AE29").Formula = '=$AC$1/AD29'
$ws.Range("AG29").Value = 857
$ws.Range("AH29").Formula = "=AG29/1000"
$ws.Range("AI29").Formula = '=$AG$1/AH29'
$ws.Range("AK29").Value = 975
$ws.Range("AL29").Formula = "=AK29/1000"
$ws.Range("AM29").Formula = '=$AK$1/AL29'
$ws.Range("AO29").Value = 1084
$ws.Range("AP29").Formula = "=AO29/1000"
$ws.Range("AQ29").Formula = '=$AO$1/AP29'

$ws.Range("AC30").Value = 768
$ws.Range("AD30").Formula = "=AC30/1000"
$ws.Range("AE30").Formula = '=$AC$1/AD30'
$ws.Range("AG30").Value = 869
$ws.Range("AH30").Formula = "=AG30/1000"
$ws.Range("AI30").Formula = '=$AG$1/AH30'
$ws.Range("AK30").Value = 982
$ws.Range("AL30").Formula = "=AK30/1000"
$ws.Range("AM30").Formula = '=$AK$1/AL30'
$ws.Range("AO30").Value = 1088
$ws.Range("AP30").Formula = "=AO30/1000"
$ws.Range("AQ30").Formula = '=$AO$1/AP30'

$ws.Range("AC31").Value = 752
$ws.Range("AD31").Formula = "=AC31/1000"
$ws.Range("AE31").Formula = '=$AC$1/AD31'
$ws.Range("AG31").Value = 862
$ws.Range("AH31").Formula = "=AG31/1000"
$ws.Range("AI31").Formula = '=$AG$1/AH31'
$ws.Range("AK31").Value = 984
$ws.Range("AL31").Formula = "=AK31/1000"
$ws.Range("AM31").Formula = '=$AK$1/AL31'
$ws.Range("AO31").Value = 1084
$ws.Range("AP31").Formula = "=AO31/1000"
$ws.Range("AQ31").Formula = '=$AO$1/AP31'

$ws.Range("AC32").Value = 744
$ws.Range("AD32").Formula = "=AC32/1000"
$ws.Range("AE32").Formula = '=$AC$1/AD32'
$ws.Range("AG32").Value = 867
$ws.Range("AH32").Formula = "=AG32/1000"
$ws.Range("AI32").Formula = '=$AG$1/AH32'
$ws.Range("AK32").Value = 978
$ws.Range("AL32").Formula = "=AK32/1000"
$ws.Range("AM32").Formula = '=$AK$1/AL32'
$ws.Range("AO32").Value = 1077
$ws.Range("AP32").Formula = "=AO32/1000"
$ws.Range("AQ32").Formula = '=$AO$1/AP32'

$ws.Range("AC33").Value = 748
$ws.Range("AD33").Formula = "=AC33/1000"
$ws.Range("AE33").Formula = '=$AC$1/AD33'
$ws.Range("AG33").Value = 852
$ws.Range("AH33").Formula = "=AG33/1000"
$ws.Range("AI33").Formula = '=$AG$1/AH33'
$ws.Range("AK33").Value = 969
$ws.Range("AL33").Formula = "=AK33/1000"
$ws.Range("AM33").Formula = '=$AK$1/AL33'
$ws.Range("AO33").Value = 1076
$ws.Range("AP33").Formula = "=AO33/1000"
$ws.Range("AQ33").Formula = '=$AO$1/AP33'

# --- Row 35: averages for new blocks ---
$ws.Range("AE35").Formula = "=AVERAGE(AE4:AE33)"
$ws.Range("AI35").Formula = "=AVERAGE(AI4:AI33)"
$ws.Range("AM35").Formula = "=AVERAGE(AM4:AM33)"
$ws.Range("AQ35").Formula = "=AVERAGE(AQ4:AQ33)"

# --- Row 36: stdev for new blocks ---
$ws.Range("AE36").Formula = "=STDEV(AE4:AE33)"
$ws.Range("AI36").Formula = "=STDEV(AI4:AI33)"
$ws.Range("AM36").Formula = "=STDEV(AM4:AM33)"
$ws.Range("AQ36").Formula = "=STDEV(AQ4:AQ33)"

# --- Row 40/41: overall mean/stdev across all blocks (now includes new ones) ---
$ws.Range("D40").Formula = "=AVERAGE(G35,K35,O35,S35,W35,AA35,AE35,AI35,AM35,AQ35)"
$ws.Range("D41").Formula = "=AVERAGE(G36,K36,O36,S36,W36,AA36,AE36,AI36,AM36,AQ36)"

Write-Host "edit complete"